$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each assignment below reproduces the exact text content from the target
# OOXML diff. Numeric-looking strings in column D (e.g. "1.004", "0.4300")
# are prefixed with a literal leading apostrophe so Excel stores them as
# TEXT (matching the original inlineStr cells) instead of re-parsing them
# as numbers (which would silently drop trailing zeros / use sci notation).
# ClearFormats() afterwards drops the transient "quote prefix" style Excel
# applies when it converts a numeric-looking entry to text, so the cell
# keeps using the default/general style exactly like the rest of the sheet.

# Row 2
$ws.Range("D2").Value = '27.395.68'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = '1.829.60'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.61%  '
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").Value = '''1.004'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -3.20%  '
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = '''315.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.59%  '
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = '''1.003'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.90%  '
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").Value = '''0.4300'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.60%  '
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = '''0.3703'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.04%  '
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").Value = '''0.07257'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").Value = '''0.8669'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").Value = '''21.18'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.08%  '
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = '1.833.55'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").Value = '''6.683'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = '''5.358'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = '''0.07100'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = '''87.74'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = '''1.006'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.28%  '
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = '''0.000008903'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.33%  '
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("D19").Value = '''1.004'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.89%  '
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("D20").Value = '''15.25'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.22%  '
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = '27.405.02'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").Value = '''10.85'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.79%  '
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = '2.057.77'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = '''2.013'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.70%  '
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").Value = '''153.25'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.86%  '
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").Value = '''18.47'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").Value = '''2.147'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +6.88%  '
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").Value = '''5.300'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").Value = '''117.42'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.24%  '
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").Value = '''0.08841'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C32").ClearFormats()
$ws.Range("D32").Value = '''0.7678'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = '''1.203'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").Value = '''4.505'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("D35").Value = '''2.879'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.81%  '
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").Value = '''1.121'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").Value = '''0.01962'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = '''0.05273'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C40").ClearFormats()
$ws.Range("D40").Value = '''2.878'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("B41").ClearFormats()
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C41").ClearFormats()
$ws.Range("D41").Value = '''7.122'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = '''0.1680'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").Value = '''0.5072'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.79%  '
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = '''8.672'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = '''10.61'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").Value = '''106.44'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.24%  '
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").Value = '''0.4733'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").Value = '''0.06425'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").Value = '''1.671'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("D51").Value = '''1.823'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.39%  '
$ws.Range("E51").ClearFormats()
